# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker "RICARDO QUINTERO SIERRA" (CC 79566378, periodo 1611) used to be
# the last data row (row 50) of the table; it now becomes the first data row
# (row 16). "LADY PAOLA VARGAS DE ORO" (CC 45528490) keeps rows 17-50 but her
# mora periods are now listed in ascending order (1812 -> 2109) instead of the
# previous descending order (2109 -> 1812).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 -> Ricardo Quintero Sierra, periodo 1611
$ws.Range("B16").Value2 = "CC"
$ws.Range("C16").Value2 = "79566378"
$ws.Range("D16").Value2 = "RICARDO QUINTERO SIERRA"
$ws.Range("E16").Value2 = "1611"
$ws.Range("F16").Value2 = 5145
$ws.Range("G16").Value2 = 3858750

# Rows 17-50 -> Lady Paola Vargas de Oro, periodos 1812..2109 ascending
$periods = @("1812","1901","1902","1903","1904","1905","1906","1907","1908","1909", `
             "1910","1911","1912","2001","2002","2003","2004","2005","2006","2007", `
             "2008","2009","2010","2011","2012","2101","2102","2103","2104","2105", `
             "2106","2107","2108","2109")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 17 + $i
    $ws.Range("B$row").Value2 = "CC"
    $ws.Range("C$row").Value2 = "45528490"
    $ws.Range("D$row").Value2 = "LADY PAOLA VARGAS DE ORO"
    $ws.Range("E$row").Value2 = $periods[$i]
    $ws.Range("F$row").Value2 = 165155
    $ws.Range("G$row").Value2 = 4760950
}
